$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H64").Value = 4389.067
$ws.Range("I64").Value = 4000
$ws.Range("K64").Value = 4000
$ws.Range("M64").Value = -3752
$ws.Range("H67").Value = 4389.067
$ws.Range("I67").Value = 4000
$ws.Range("K67").Value = 4000
$ws.Range("M67").Value = -3142
$ws.Range("H76").Value = 3442.8572
$ws.Range("I76").Value = 3466.6667
$ws.Range("J76").Value = 3300
$ws.Range("K76").Value = 3466.6667
$ws.Range("L76").Value = 3300
$ws.Range("M76").Value = -3151.6667
$ws.Range("N76").Value = -3930
$ws.Range("H79").Value = 3442.8572
$ws.Range("I79").Value = 3466.6667
$ws.Range("J79").Value = 3300
$ws.Range("K79").Value = 3466.6667
$ws.Range("L79").Value = 3300
$ws.Range("M79").Value = -2374.6667
$ws.Range("N79").Value = -5484
$ws.Range("H116").Value = 3101.5
$ws.Range("I116").Value = 2080
$ws.Range("J116").Value = 4804
$ws.Range("K116").Value = 2080
$ws.Range("L116").Value = 4804
$ws.Range("M116").Value = 1362
$ws.Range("N116").Value = -11688
$ws.Range("H138").Value = 11115454
$ws.Range("I138").Value = 8898.200000000001
$ws.Range("J138").Value = 15387207
$ws.Range("K138").Value = 26694.6
$ws.Range("L138").Value = 46161621
$ws.Range("M138").Value = -21554.6
$ws.Range("N138").Value = -46171901

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3013
$ws.Range("I63").Value = 3054.0908
$ws.Range("J63").Value = 2900
$ws.Range("K63").Value = 3054.0908
$ws.Range("L63").Value = 2900
$ws.Range("M63").Value = -2368.0908
$ws.Range("N63").Value = -4272
$ws.Range("H66").Value = 3013
$ws.Range("I66").Value = 3054.0908
$ws.Range("J66").Value = 2900
$ws.Range("K66").Value = 15270.454
$ws.Range("L66").Value = 14500
$ws.Range("M66").Value = -11838.454
$ws.Range("N66").Value = -21364
$ws.Range("H110").Value = 1106.2963
$ws.Range("I110").Value = 970.6316
$ws.Range("J110").Value = 1428.5
$ws.Range("K110").Value = 970.6316
$ws.Range("L110").Value = 1428.5
$ws.Range("M110").Value = 1074.3684
$ws.Range("N110").Value = -5518.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1816.5577
$ws.Range("I134").Value = 1732.289
$ws.Range("J134").Value = 2358.2856
$ws.Range("K134").Value = 5196.867
$ws.Range("L134").Value = 7074.8568
$ws.Range("M134").Value = -2661.867
$ws.Range("N134").Value = -12144.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 633.6829
$ws.Range("I5").Value = 265.65625
$ws.Range("J5").Value = 1942.2222
$ws.Range("K5").Value = 796.96875
$ws.Range("L5").Value = 5826.6666
$ws.Range("M5").Value = -684.96875
$ws.Range("N5").Value = -6050.6666
$ws.Range("H122").Value = 1074.2609
$ws.Range("I122").Value = 343.7143
$ws.Range("K122").Value = 3093.4287
$ws.Range("M122").Value = -643.4286999999999
$ws.Range("H131").Value = 984.1667
$ws.Range("J131").Value = 1025.8928
$ws.Range("L131").Value = 3077.6784
$ws.Range("N131").Value = -13157.6784
$ws.Range("H135").Value = 633.6829
$ws.Range("I135").Value = 265.65625
$ws.Range("J135").Value = 1942.2222
$ws.Range("K135").Value = 2390.90625
$ws.Range("L135").Value = 17479.9998
$ws.Range("M135").Value = 144.09375
$ws.Range("N135").Value = -22549.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H10").Value = 5000300
$ws.Range("I10").Value = 5000300
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 5000300
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -5000131
$ws.Range("N10").ClearContents()
$ws.Range("H70").Value = 83677.62
$ws.Range("I70").Value = 172650
$ws.Range("J70").Value = 7415.5713
$ws.Range("K70").Value = 172650
$ws.Range("L70").Value = 7415.5713
$ws.Range("M70").Value = -172380
$ws.Range("N70").Value = -7955.5713
$ws.Range("H73").Value = 83677.62
$ws.Range("I73").Value = 172650
$ws.Range("J73").Value = 7415.5713
$ws.Range("K73").Value = 172650
$ws.Range("L73").Value = 7415.5713
$ws.Range("M73").Value = -171714
$ws.Range("N73").Value = -9287.5713
$ws.Range("H122").Value = 1442.0667
$ws.Range("I122").Value = 1409.0714
$ws.Range("J122").Value = 1904
$ws.Range("K122").Value = 4227.2142
$ws.Range("L122").Value = 5712
$ws.Range("M122").Value = -1777.2142
$ws.Range("N122").Value = -10612
$ws.Range("H123").Value = 37420.727
$ws.Range("J123").Value = 37420.727
$ws.Range("L123").Value = 37420.727
$ws.Range("N123").Value = -42320.727

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2992.4167
$ws.Range("I40").Value = 2991.2727
$ws.Range("J40").Value = 3005
$ws.Range("K40").Value = 2991.2727
$ws.Range("L40").Value = 3005
$ws.Range("M40").Value = -2855.2727
$ws.Range("N40").Value = -3277
$ws.Range("H61").Value = 1627.9
$ws.Range("I61").Value = 1418
$ws.Range("J61").Value = 2467.5
$ws.Range("K61").Value = 1418
$ws.Range("L61").Value = 2467.5
$ws.Range("M61").Value = -1216
$ws.Range("N61").Value = -2871.5
$ws.Range("H93").Value = 918.4
$ws.Range("I93").Value = 899.4583
$ws.Range("K93").Value = 899.4583
$ws.Range("M93").Value = 348.5417
$ws.Range("H113").Value = 1627.9
$ws.Range("I113").Value = 1418
$ws.Range("J113").Value = 2467.5
$ws.Range("K113").Value = 1418
$ws.Range("L113").Value = 2467.5
$ws.Range("M113").Value = 752
$ws.Range("N113").Value = -6807.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 779.4583
$ws.Range("I113").Value = 773.0454999999999
$ws.Range("J113").Value = 850
$ws.Range("K113").Value = 2319.1365
$ws.Range("L113").Value = 2550
$ws.Range("M113").Value = -149.1364999999996
$ws.Range("N113").Value = -6890
$ws.Range("H126").Value = 1196.0834
$ws.Range("I126").Value = 938.6667
$ws.Range("J126").Value = 1968.3334
$ws.Range("K126").Value = 2816.0001
$ws.Range("L126").Value = 5905.0002
$ws.Range("M126").Value = -346.0001000000002
$ws.Range("N126").Value = -10845.0002
